# Generate Report for Handback
# - Updates the "Status" text from "Ready for handoff" to
#   "Handed back: in sync with en-US" everywhere it appears.
# - Stamps the "Latest Handback DateTime" column (H) with real timestamps
#   (previously the zero-date placeholder) per language sheet.
# - Adds the "Latest Target File" (F) and "Latest Handback File" (G)
#   columns with hyperlinked file names for each data row, on both the
#   zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/0cd3d88d511f61ad900b0511f14ce547934f1055/e2e/59db17ff-0d89-493f-a4fb-64bdf414a197.md"
$mdDisplay = "59db17ff-0d89-493f-a4fb-64bdf414a197.md"

# ---------------------------------------------------------------------
# 1. Overview sheet: roll the new status text into the summary columns.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Latest Handback DateTime column now has a real timestamp
$wsZh.Range("H2").Value = "2016-03-22 17:13:23"
$wsZh.Range("H3").Value = "2016-03-22 17:13:23"

$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfc61ff083b488bac160d5f280d22a652d52b7b0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.zh-cn.xlf"
$zhXlfDisplay = "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.zh-cn.xlf"

# Latest Target File / Latest Handback File columns, row 2
$wsZh.Range("F2").Value = $mdDisplay
$wsZh.Range("F2").Font.Underline = 2
$wsZh.Range("F2").Font.Color = 15570276
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $mdAddress, "", "", $mdDisplay)

$wsZh.Range("G2").Value = $zhXlfDisplay
$wsZh.Range("G2").Font.Underline = 2
$wsZh.Range("G2").Font.Color = 15570276
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfAddress, "", "", $zhXlfDisplay)

# Latest Target File / Latest Handback File columns, row 3
# (row 3's source file includes/depends on row 2's file, so it reuses the
# same target/handback artifacts as row 2)
$wsZh.Range("F3").Value = $mdDisplay
$wsZh.Range("F3").Font.Underline = 2
$wsZh.Range("F3").Font.Color = 15570276
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $mdAddress, "", "", $mdDisplay)

$wsZh.Range("G3").Value = $zhXlfDisplay
$wsZh.Range("G3").Font.Underline = 2
$wsZh.Range("G3").Font.Color = 15570276
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfAddress, "", "", $zhXlfDisplay)

# ---------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Latest Handback DateTime column now has a real timestamp
$wsDe.Range("H2").Value = "2016-03-22 17:13:30"
$wsDe.Range("H3").Value = "2016-03-22 17:13:30"

$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/83edba7cfcf10a304a23586e26fbfa94e4a18fcb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.de-de.xlf"
$deXlfDisplay = "59db17ff-0d89-493f-a4fb-64bdf414a197.bac44bfe36bffe9cc476143af072f2ebdc47248a.de-de.xlf"

# Latest Target File / Latest Handback File columns, row 2
$wsDe.Range("F2").Value = $mdDisplay
$wsDe.Range("F2").Font.Underline = 2
$wsDe.Range("F2").Font.Color = 15570276
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $mdAddress, "", "", $mdDisplay)

$wsDe.Range("G2").Value = $deXlfDisplay
$wsDe.Range("G2").Font.Underline = 2
$wsDe.Range("G2").Font.Color = 15570276
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfAddress, "", "", $deXlfDisplay)

# Latest Target File / Latest Handback File columns, row 3
$wsDe.Range("F3").Value = $mdDisplay
$wsDe.Range("F3").Font.Underline = 2
$wsDe.Range("F3").Font.Color = 15570276
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $mdAddress, "", "", $mdDisplay)

$wsDe.Range("G3").Value = $deXlfDisplay
$wsDe.Range("G3").Font.Underline = 2
$wsDe.Range("G3").Font.Color = 15570276
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfAddress, "", "", $deXlfDisplay)
